$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = "23.127.75"
$ws.Range("E2").Value = "  -3.66%  "

# Row 3
$ws.Range("D3").Value = "1.602.68"
$ws.Range("E3").Value = "  -2.97%  "

# Row 4
$ws.Range("E4").Value = "  +0.07%  "

# Row 5
$ws.Range("E5").Value = "  +0.01%  "

# Row 6
$ws.Range("D6").Value = "'301.42"
$ws.Range("E6").Value = "  -2.80%  "

# Row 7
$ws.Range("E7").Value = "  -3.37%  "

# Row 8
$ws.Range("D8").Value = "'0.3653"
$ws.Range("E8").Value = "  -4.26%  "

# Row 9
$ws.Range("D9").Value = "'50.03"
$ws.Range("E9").Value = "  -4.46%  "

# Row 10
$ws.Range("D10").Value = "'1.257"
$ws.Range("E10").Value = "  -6.89%  "

# Row 11
$ws.Range("E11").Value = "  +0.10%  "

# Row 12
$ws.Range("D12").Value = "'0.08139"
$ws.Range("E12").Value = "  -3.76%  "

# Row 13
$ws.Range("E13").Value = "  -3.74%  "

# Row 14
$ws.Range("D14").Value = "'6.596"
$ws.Range("E14").Value = "  -6.75%  "

# Row 15
$ws.Range("B15").Value = "ShibaInu"
$ws.Range("C15").Value = "https://coinranking.com/coin/xz24e0BjL+shibainu-shib"
$ws.Range("D15").Value = "'0.00001257"
$ws.Range("E15").Value = "  -4.13%  "

# Row 16
$ws.Range("B16").Value = "Chainlink"
$ws.Range("C16").Value = "https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link"
$ws.Range("D16").Value = "'7.408"
$ws.Range("E16").Value = "  -7.58%  "

# Row 17
$ws.Range("D17").Value = "1.603.65"
$ws.Range("E17").Value = "  -2.82%  "

# Row 18
$ws.Range("D18").Value = "'91.61"

# Row 19
$ws.Range("D19").Value = "'0.06852"
$ws.Range("E19").Value = "  -2.24%  "

# Row 20
$ws.Range("D20").Value = "'18.30"
$ws.Range("E20").Value = "  -7.28%  "

# Row 21
$ws.Range("E21").Value = "  -6.08%  "

# Row 22
$ws.Range("B22").Value = "Dai"
$ws.Range("C22").Value = "https://coinranking.com/coin/MoTuySvg7+dai-dai"
$ws.Range("D22").Value = "'1.001"
$ws.Range("E22").Value = "  -0.04%  "

# Row 23
$ws.Range("B23").Value = "Cosmos"
$ws.Range("C23").Value = "https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom"
$ws.Range("D23").Value = "'13.00"
$ws.Range("E23").Value = "  -5.84%  "

# Row 24
$ws.Range("B24").Value = "WrappedBTC"
$ws.Range("C24").Value = "https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc"
$ws.Range("D24").Value = "23.139.33"
$ws.Range("E24").Value = "  -3.55%  "

# Row 25
$ws.Range("B25").Value = "Toncoin"
$ws.Range("C25").Value = "https://coinranking.com/coin/67YlI0K1b+toncoin-ton"
$ws.Range("D25").Value = "'2.337"
$ws.Range("E25").Value = "  -4.34%  "

# Row 26
$ws.Range("B26").Value = "LidoDAOToken"
$ws.Range("C26").Value = "https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo"
$ws.Range("D26").Value = "'2.730"
$ws.Range("E26").Value = "  -7.56%  "

# Row 27
$ws.Range("B27").Value = "EthereumClassic"
$ws.Range("C27").Value = "https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc"
$ws.Range("D27").Value = "'21.10"
$ws.Range("E27").Value = "  -4.53%  "

# Row 28
$ws.Range("B28").Value = "Monero"
$ws.Range("C28").Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$ws.Range("D28").Value = "'150.26"
$ws.Range("E28").Value = "  -1.75%  "

# Row 29
$ws.Range("B29").Value = "HuobiToken"
$ws.Range("C29").Value = "https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht"
$ws.Range("D29").Value = "'5.290"
$ws.Range("E29").Value = "  -2.19%  "

# Row 30
$ws.Range("B30").Value = "BitcoinCash"
$ws.Range("C30").Value = "https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch"
$ws.Range("D30").Value = "'132.02"
$ws.Range("E30").Value = "  -4.47%  "

# Row 31
$ws.Range("B31").Value = "WEMIXTOKEN"
$ws.Range("C31").Value = "https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix"
$ws.Range("D31").Value = "'2.431"
$ws.Range("E31").Value = "  -3.02%  "

# Row 32
$ws.Range("B32").Value = "Filecoin"
$ws.Range("C32").Value = "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
$ws.Range("D32").Value = "'6.823"
$ws.Range("E32").Value = "  -14.07%  "

# Row 33
$ws.Range("B33").Value = "WrappedliquidstakedEther2.0"
$ws.Range("C33").Value = "https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth"
$ws.Range("D33").Value = "1.780.76"
$ws.Range("E33").Value = "  -2.74%  "

# Row 34
$ws.Range("B34").Value = "Hedera"
$ws.Range("C34").Value = "https://coinranking.com/coin/jad286TjB+hedera-hbar"
$ws.Range("D34").Value = "'0.07704"
$ws.Range("E34").Value = "  -4.42%  "

# Row 35
$ws.Range("B35").Value = "ImmutableX"
$ws.Range("C35").Value = "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
$ws.Range("D35").Value = "'0.9497"
$ws.Range("E35").Value = "  -6.98%  "

# Row 36
$ws.Range("B36").Value = "VeChain"
$ws.Range("C36").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D36").Value = "'0.02745"
$ws.Range("E36").Value = "  -6.20%  "

# Row 37
$ws.Range("B37").Value = "InternetComputer(DFINITY)"
$ws.Range("C37").Value = "https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp"
$ws.Range("D37").Value = "'6.259"
$ws.Range("E37").Value = "  -7.33%  "

# Row 38
$ws.Range("B38").Value = "Algorand"
$ws.Range("C38").Value = "https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo"
$ws.Range("D38").Value = "'0.2546"
$ws.Range("E38").Value = "  -5.00%  "

# Row 39
$ws.Range("B39").Value = "Stellar"
$ws.Range("C39").Value = "https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"
$ws.Range("D39").Value = "'0.08909"
$ws.Range("E39").Value = "  -2.12%  "

# Row 40
$ws.Range("B40").Value = "FraxShare"
$ws.Range("C40").Value = "https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs"
$ws.Range("D40").Value = "'10.05"
$ws.Range("E40").Value = "  -6.57%  "

# Row 41
$ws.Range("B41").Value = "TrustWalletToken"
$ws.Range("C41").Value = "https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt"
$ws.Range("D41").Value = "'1.387"
$ws.Range("E41").Value = "  -2.55%  "

# Row 42
$ws.Range("B42").Value = "TheSandbox"
$ws.Range("C42").Value = "https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand"
$ws.Range("D42").Value = "'0.7103"
$ws.Range("E42").Value = "  -6.78%  "

# Row 43
$ws.Range("B43").Value = "Aptos"
$ws.Range("C43").Value = "https://coinranking.com/coin/HGYj5JCv5+aptos-apt"
$ws.Range("D43").Value = "'12.70"
$ws.Range("E43").Value = "  -5.42%  "

# Row 44
$ws.Range("B44").Value = "EnergySwap"
$ws.Range("C44").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D44").Value = "'15.53"
$ws.Range("E44").Value = "  -4.85%  "

# Row 45
$ws.Range("B45").Value = "Decentraland"
$ws.Range("C45").Value = "https://coinranking.com/coin/tEf7-dnwV3BXS+decentraland-mana"
$ws.Range("D45").Value = "'0.6628"
$ws.Range("E45").Value = "  -5.04%  "

# Row 46
$ws.Range("B46").Value = "Frax"
$ws.Range("C46").Value = "https://coinranking.com/coin/KfWtaeV1W+frax-frax"
$ws.Range("D46").Value = "'1.0000"
$ws.Range("E46").Value = "  +0.02%  "

# Row 47
$ws.Range("B47").Value = "NEARProtocol"
$ws.Range("C47").Value = "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
$ws.Range("D47").Value = "'2.304"
$ws.Range("E47").Value = "  -6.26%  "

# Row 48
$ws.Range("B48").Value = "PancakeSwap"
$ws.Range("C48").Value = "https://coinranking.com/coin/ncYFcP709+pancakeswap-cake"
$ws.Range("D48").Value = "'3.977"
$ws.Range("E48").Value = "  -3.06%  "

# Row 49
$ws.Range("B49").Value = "Quant"
$ws.Range("C49").Value = "https://coinranking.com/coin/bauj_21eYVwso+quant-qnt"
$ws.Range("D49").Value = "'132.07"
$ws.Range("E49").Value = "  -2.16%  "

# Row 50
$ws.Range("B50").Value = "Cronos"
$ws.Range("C50").Value = "https://coinranking.com/coin/65PHZTpmE55b+cronos-cro"
$ws.Range("D50").Value = "'0.07949"
$ws.Range("E50").Value = "  -4.46%  "

# Row 51
$ws.Range("B51").Value = "Flow"
$ws.Range("C51").Value = "https://coinranking.com/coin/QQ0NCmjVq+flow-flow"
$ws.Range("D51").Value = "'1.209"
$ws.Range("E51").Value = "  -2.03%  "
